$d = $word.ActiveDocument

# 1) Fix the "computer functions" text - remove the bookmark split and join the text
$d.Content.Find.Execute("Supported faculty in basic computer functions required to conduct lectures in high-tech classroom environments with state-of-the-art audio/video equipment", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

